# Updated cryptos list data (prices + 1h volume deltas) pulled from coinranking.com
# Applies the latest snapshot values to Sheet1, preserving text-formatted cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.727.79"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "2.202.06"
$ws.Range("E3").Value = "  -2.77%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.82"
$ws.Range("E5").Value = "  -1.93%  "

$ws.Range("E6").Value = "  -4.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.18"
$ws.Range("E7").Value = "  -5.53%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -2.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.05"
$ws.Range("E10").Value = "  -4.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.103"
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").Value = "2.533.04"
$ws.Range("E13").Value = "  -2.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.34"
$ws.Range("E14").Value = "  -4.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.15"
$ws.Range("E15").Value = "  -3.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.56"
$ws.Range("E16").Value = "  -2.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  -3.69%  "

$ws.Range("D18").Value = "2.198.40"
$ws.Range("E18").Value = "  -2.88%  "

$ws.Range("D19").Value = "41.627.16"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("E20").Value = "  -3.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.04"
$ws.Range("E21").Value = "  -3.81%  "

$ws.Range("E22").Value = "  -2.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.48"
$ws.Range("E23").Value = "  -4.58%  "

$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("E25").Value = "  -3.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  -3.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  -1.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.56"
$ws.Range("E28").Value = "  -1.40%  "

$ws.Range("E29").Value = "  -6.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.45"
$ws.Range("E30").Value = "  -0.47%  "

$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.61"
$ws.Range("E32").Value = "  -8.10%  "

$ws.Range("E33").Value = "  -3.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.98"
$ws.Range("E34").Value = "  -2.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.60"
$ws.Range("E35").Value = "  -4.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0644"
$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("E37").Value = "  -7.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.28"
$ws.Range("E38").Value = "  -8.68%  "

$ws.Range("E39").Value = "  -4.66%  "

$ws.Range("B40").Value = "TerraClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000238"
$ws.Range("E40").Value = "  -9.72%  "

$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("E42").Value = "  -1.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.64"
$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.43"
$ws.Range("E44").Value = "  -13.93%  "

$ws.Range("E45").Value = "  -2.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  -3.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "96.76"
$ws.Range("E47").Value = "  -5.30%  "

$ws.Range("D48").Value = "1.460.80"
$ws.Range("E48").Value = "  -2.99%  "

$ws.Range("E49").Value = "  -1.83%  "

$ws.Range("E50").Value = "  -9.34%  "

$ws.Range("E51").Value = "  -5.54%  "

